# Updated SAD with container pattern and Risks
# Adjust the "Risk Probability of Occurrence" values in the Risk evaluation
# table; the dependent "Risk Factor" formulas (E7:E11 = C*D) recalculate
# automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Risk evaluation")

$ws.Range("C7").Value  = 0.16
$ws.Range("C8").Value  = 0.3
$ws.Range("C9").Value  = 0.12
$ws.Range("C10").Value = 0.2
$ws.Range("C11").Value = 0.01

# Move the view / active selection to reflect where the author was working.
$excel.ActiveWindow.ScrollRow    = 6
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("C6").Select()
